# Commit: "Add files via upload"
# The uploaded workbook replaces the father's-phone ("رقم الأب", column C)
# values for students in rows 16-19 with a single shared number, and
# replaces the teacher name ("اسم المدرس", column F) throughout the sheet
# from the Arabic "عبد الرحمن سعيد" to "MR. Abdulrahman Saeed". Also moves
# the active selection to G5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("رقم الأب") for rows 16-19 -> unified phone number
$ws.Range("C16").Value = "011 5101 4252"
$ws.Range("C17").Value = "011 5101 4252"
$ws.Range("C18").Value = "011 5101 4252"
$ws.Range("C19").Value = "011 5101 4252"

# Column F ("اسم المدرس") for every data row -> new teacher name
$ws.Range("F2:F27").Value = "MR. Abdulrahman Saeed"

# Move the active selection, as reflected in the saved sheetView
$null = $ws.Range("G5").Select()
